# Update "想去人数" (want-to-go count) figures in column F across sheets,
# as published in the refreshed gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    3  = 194
    4  = 49
    5  = 1619
    6  = 3229
    7  = 762
    8  = 1975
    9  = 1894
    10 = 970
    11 = 339
    12 = 13
    13 = 1570
    14 = 334
    17 = 1386
    18 = 479
    19 = 592
    20 = 282
    21 = 10425
    22 = 9605
    23 = 830
    24 = 633
    25 = 1802
    26 = 132
    27 = 363
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 0

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    4  = 194
    6  = 49
    7  = 1619
    8  = 3229
    9  = 762
    10 = 1975
    11 = 1894
    12 = 970
    13 = 339
    14 = 13
    15 = 1570
    16 = 334
    21 = 1386
    22 = 479
    23 = 592
    24 = 282
    25 = 10425
    26 = 9605
    27 = 830
    28 = 633
    29 = 1802
    32 = 132
    33 = 363
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
